# Security vulnerability check rows: append a new row 40 (matching the
# existing column layout: time / total-len / ID / actual-len / checksum
# and their _DEC numeric counterparts) to each of the four sensor sheets.
#
# Columns A-E, G hold text (raw hex strings / IDs / timestamps stored as
# literal text, not as Excel dates/numbers) while F, H, I are plain
# numbers - matching the pattern already used by every other data row in
# these sheets. To stop Excel's automatic type-inference from turning the
# long numeric-looking strings (column G) into floating point numbers, the
# text columns are briefly forced to Text format ("@") before the value is
# assigned, then restored to the default "Normal" style so the written
# cells stay plain (no leftover style index), exactly like the neighbouring
# rows.
#
# NOTE: named parameter binding (-Foo bar) is unreliable in this COM-interop
# runtime, so the helper function below takes plain positional parameters.

function Add-SensorRow($ws, $Row, $A, $B, $C, $D, $E, $F, $G, $H, $I) {
    $textCols = @("A", "B", "C", "D", "E", "G")
    foreach ($col in $textCols) {
        $ws.Range("$col$Row").NumberFormat = "@"
    }

    $ws.Range("A$Row").Value = $A
    $ws.Range("B$Row").Value = $B
    $ws.Range("C$Row").Value = $C
    $ws.Range("D$Row").Value = $D
    $ws.Range("E$Row").Value = $E
    $ws.Range("F$Row").Value = $F
    $ws.Range("G$Row").Value = $G
    $ws.Range("H$Row").Value = $H
    $ws.Range("I$Row").Value = $I

    foreach ($col in $textCols) {
        $ws.Range("$col$Row").Style = "Normal"
    }
}

$wb = $excel.ActiveWorkbook

# ROW35-FE-LIFTER
$ws1 = $wb.Worksheets.Item(1)
Add-SensorRow $ws1 40 "2025-03-05 23:42:06" "0x01,0x90 " "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c," "0x01,0x90," "0x d" 400 "568631262647113770877196" 400 13

# ROW35-MID-LIFTER
$ws2 = $wb.Worksheets.Item(2)
Add-SensorRow $ws2 40 "2025-03-05 23:29:35" "0x01,0x90 " "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c," "0x01,0x90," "0x e" 400 "568631262647113770942732" 400 14

# ROW02-FE-LIFTER
$ws3 = $wb.Worksheets.Item(3)
Add-SensorRow $ws3 40 "2025-03-05 23:51:45" "0x01,0x90 " "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x06,0x41,0x0c," "0x01,0x90," "0xff" 400 "568631262647113769959692" 400 255

# ROW02-MID-LIFTER
$ws4 = $wb.Worksheets.Item(4)
Add-SensorRow $ws4 40 "2025-03-05 23:41:15" "0x01,0x90 " "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c," "0x01,0x90," "0x 3" 400 "568631262647113769959692" 400 3
